$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Tiny precision correction on H13 (floating point re-computation artifact)
$ws.Range("H13").Value = 0.9890787327105353

# New row 16: 14th "HexGrid-60degTilt5degRes" entry from the Gaussian
# Quadrature scheme export
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.781294520948619
$ws.Range("D16").Value = 1.308354768009678
$ws.Range("E16").Value = 0.9245862354206003
$ws.Range("F16").Value = 1.781294520948619
$ws.Range("G16").Value = 0.9667450792716777
$ws.Range("H16").Value = 1.036695227981307
$ws.Range("I16").Value = 0.9440024818372539
$ws.Range("J16").Value = 1.308354768009678
$ws.Range("K16").Value = 1.116470501715139
$ws.Range("L16").Value = 1.448882511331879
$ws.Range("M16").Value = 1.160279718911523

# Match style of the A column "index" cells (bold, bordered, centered)
# by copying formats from the cell above (keeps the value already set)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
